$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: relabel/reorder the property-type columns ---
$ws.Range("D1").Value = "Propriétées : seulement copropriétées"
$ws.Range("B1").Value = "Propriétées : toutes"
$ws.Range("C1").Value = "Propriétées : seulement individuelles"

# --- Row 3 label change ("Seulement Individuelles" -> "2 (surface + chambres)") ---
$ws.Range("A3").Value = "2 (surface + chambres)"

# --- Row 5: C5/D5 relabeled from "Both" to "N/A" ---
$ws.Range("C5").Value = "N/A"
$ws.Range("D5").Value = "N/A"

# --- New row 6: averages ---
$ws.Range("A6").Value = "Moyenne :"
$ws.Range("B6").Formula = "=AVERAGE(B2:B5)"
$ws.Range("C6").Formula = "=AVERAGE(C2:C4)"
$ws.Range("D6").Formula = "=AVERAGE(D2:D4)"

# Match formatting of the rows above (A column left-aligned, others centered)
$ws.Range("A6").HorizontalAlignment = $ws.Range("A5").HorizontalAlignment
$ws.Range("B6").HorizontalAlignment = $ws.Range("B5").HorizontalAlignment
$ws.Range("C6").HorizontalAlignment = $ws.Range("C5").HorizontalAlignment
$ws.Range("D6").HorizontalAlignment = $ws.Range("D5").HorizontalAlignment

# --- Column widths (best-fit to new, longer header text) ---
$ws.Columns.Item(2).ColumnWidth = 18.333333333333332
$ws.Columns.Item(3).ColumnWidth = 34.5
$ws.Columns.Item(4).ColumnWidth = 34.166666666666664

# --- Selection moved ---
$ws.Range("A11").Select()
